# Regenerate orders with updated distance/size codes.
# Distance codes: D80->D86, D51->D55, D64->D69
# Size code:      S30->S31
# These substitutions occur inside composite strings such as
# "Face05_D80_S25", "Fixation_D80_l.png", the standalone "D80" / "S30"
# lookup values, etc. A simple substring Find & Replace across every used
# cell on the sheet reproduces the diff exactly, since the cell layout,
# row/column structure and all other text are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = $ws.UsedRange

$cells.Replace("D80", "D86")
$cells.Replace("D51", "D55")
$cells.Replace("D64", "D69")
$cells.Replace("S30", "S31")
